# Add a new column "mX0" (derived from cX0) right after the existing
# "cX0" column (column H) in the metadata sheet.
#
# This mirrors the author's edit: select/insert a blank column at I,
# which shifts every column from I..S to J..T, then populate the new
# column with a header, a unit label, and a formula (=cX0*0.5) for
# each data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column before the current column I (cS0).
# Excel shifts I:S -> J:T and copies column-H's formatting onto the
# newly inserted column.
$ws.Columns("I:I").Insert()

# Match column H's width on the freshly inserted column I (Excel does
# this automatically when inserting; make it explicit here).
$ws.Columns("I:I").ColumnWidth = $ws.Columns("H:H").ColumnWidth()

# Row 1 header for the new column.
$ws.Range("I1").Value = "mX0"

# Row 2 holds the unit for each column; mX0 is a mass, same unit ("g")
# as the neighbouring mS0 column (now K2).
$ws.Range("I2").Value = $ws.Range("K2").Value()

# Data rows: mX0 = cX0 * 0.5 (same 0.5 factor used by the existing
# mS0 = cS0 * 0.5 formula). Row 3 is entered on its own; rows 4-7
# share one formula. The sibling mS0 formulas in column K (shifted
# from the old J) are re-entered the same way so they keep referring
# to the correct (shifted) cS0 column (J).
$ws.Range("I3").Formula = "=H3*0.5"
$ws.Range("I4:I7").Formula = "=H4*0.5"
$ws.Range("K4:K7").Formula = "=J4*0.5"

# Leave the newly inserted column selected, matching the author's
# final selection after inserting the column.
$ws.Columns("I:I").Select()
